# The "EU_estimate" column (W) and the "EU_estimate_Limit" column (AB) are
# being removed from the Databank sheet. Deleting column W first shifts
# every column to its right one place to the left (X->W, Y->X, Z->Y,
# AA->Z, AB->AA); the old "EU_estimate_Limit" data therefore now lives in
# column AA, so that is the second column to delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("W").Delete()
$ws.Columns("AA").Delete()
